$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3726.7273
$ws.Range("I40").Value = 2499.6667
$ws.Range("J40").Value = 4186.875
$ws.Range("K40").Value = 2499.6667
$ws.Range("L40").Value = 4186.875
$ws.Range("M40").Value = -2324.6667
$ws.Range("N40").Value = -4536.875

$ws.Range("H57").Value = 76326.336
$ws.Range("J57").Value = 76326.336
$ws.Range("L57").Value = 228979.008
$ws.Range("N57").Value = -229977.008

$ws.Range("H62").Value = 3199.182
$ws.Range("J62").Value = 2414.6667
$ws.Range("L62").Value = 2414.6667
$ws.Range("N62").Value = -3662.6667

$ws.Range("H65").Value = 3199.182
$ws.Range("J65").Value = 2414.6667
$ws.Range("L65").Value = 12073.3335
$ws.Range("N65").Value = -18313.3335

$ws.Range("H103").Value = 647.75
$ws.Range("I103").Value = 552
$ws.Range("J103").Value = 743.5
$ws.Range("K103").Value = 1656
$ws.Range("L103").Value = 2230.5
$ws.Range("M103").Value = -1070
$ws.Range("N103").Value = -3402.5

$ws.Range("H111").Value = 1407.3334
$ws.Range("I111").Value = 1130
$ws.Range("J111").Value = 1962
$ws.Range("K111").Value = 3390
$ws.Range("L111").Value = 5886
$ws.Range("M111").Value = -323
$ws.Range("N111").Value = -12020

$ws.Range("H132").Value = 7982.857
$ws.Range("I132").Value = 8167.185
$ws.Range("K132").Value = 24501.555
$ws.Range("M132").Value = -21971.555

$ws.Range("H138").Value = 4060.9546
$ws.Range("I138").Value = 4060.1428
$ws.Range("K138").Value = 12180.4284
$ws.Range("M138").Value = -7040.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6602.9116
$ws.Range("I74").Value = 2677
$ws.Range("J74").Value = 19362.125
$ws.Range("K74").Value = 2677
$ws.Range("L74").Value = 19362.125
$ws.Range("M74").Value = -1803
$ws.Range("N74").Value = -21110.125

$ws.Range("H77").Value = 6602.9116
$ws.Range("I77").Value = 2677
$ws.Range("J77").Value = 19362.125
$ws.Range("K77").Value = 13385
$ws.Range("L77").Value = 96810.625
$ws.Range("M77").Value = -9017
$ws.Range("N77").Value = -105546.625

$ws.Range("H97").Value = 3066.3333
$ws.Range("I97").Value = 2999.5
$ws.Range("K97").Value = 2999.5
$ws.Range("M97").Value = -2503.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 96550
$ws.Range("J59").Value = 96550
$ws.Range("L59").Value = 96550
$ws.Range("N59").Value = -98244

$ws.Range("H94").Value = 4313.4443
$ws.Range("I94").Value = 4474.4287
$ws.Range("K94").Value = 4474.4287
$ws.Range("M94").Value = -4023.4287

$ws.Range("H107").Value = 1470.1578
$ws.Range("I107").Value = 1169.9524
$ws.Range("J107").Value = 1841
$ws.Range("K107").Value = 1169.9524
$ws.Range("L107").Value = 1841
$ws.Range("M107").Value = 750.0476000000001
$ws.Range("N107").Value = -5681

$ws.Range("H134").Value = 4407.7
$ws.Range("I134").Value = 4341.8887
$ws.Range("K134").Value = 13025.6661
$ws.Range("M134").Value = -10490.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 62356
$ws.Range("J52").Value = 65445
$ws.Range("L52").Value = 65445
$ws.Range("N52").Value = -66033

$ws.Range("H99").Value = 2365.7
$ws.Range("I99").Value = 2293.1667
$ws.Range("J99").Value = 2474.5
$ws.Range("K99").Value = 2293.1667
$ws.Range("L99").Value = 2474.5
$ws.Range("M99").Value = -795.1667000000002
$ws.Range("N99").Value = -5470.5

$ws.Range("H116").Value = 59247.332
$ws.Range("J116").Value = 59247.332
$ws.Range("L116").Value = 59247.332
$ws.Range("N116").Value = -68425.33199999999

$ws.Range("H122").Value = 1997.2069
$ws.Range("J122").Value = 2499.4
$ws.Range("L122").Value = 7498.200000000001
$ws.Range("N122").Value = -12398.2

$ws.Range("H126").Value = 2365.7
$ws.Range("I126").Value = 2293.1667
$ws.Range("J126").Value = 2474.5
$ws.Range("K126").Value = 6879.500100000001
$ws.Range("L126").Value = 7423.5
$ws.Range("M126").Value = -4409.500100000001
$ws.Range("N126").Value = -12363.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1391.7222
$ws.Range("I5").Value = 1235.0834
$ws.Range("K5").Value = 3705.2502
$ws.Range("M5").Value = -3593.2502

$ws.Range("H7").Value = 71
$ws.Range("I7").Value = 71
$ws.Range("K7").Value = 213
$ws.Range("M7").Value = -101

$ws.Range("H9").Value = 330398.16
$ws.Range("J9").Value = 301547.75
$ws.Range("L9").Value = 904643.25
$ws.Range("N9").Value = -905091.25

$ws.Range("H34").Value = 2199.2354
$ws.Range("J34").Value = 2303.3572
$ws.Range("L34").Value = 6910.071599999999
$ws.Range("N34").Value = -7078.071599999999

$ws.Range("H39").Value = 14521.417
$ws.Range("J39").Value = 43424.5
$ws.Range("L39").Value = 130273.5
$ws.Range("N39").Value = -130861.5

$ws.Range("H55").Value = 53127604
$ws.Range("I55").Value = 120000616
$ws.Range("J55").Value = 1115265.9
$ws.Range("K55").Value = 360001848
$ws.Range("L55").Value = 3345797.7
$ws.Range("M55").Value = -360001671
$ws.Range("N55").Value = -3346151.7

$ws.Range("H112").Value = 8099.8335
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = 8836.091
$ws.Range("K112").Value = 3
$ws.Range("L112").Value = 26508.273
$ws.Range("M112").Value = 1105
$ws.Range("N112").Value = -28724.273

$ws.Range("H120").Value = 17562.375
$ws.Range("I120").Value = 11875
$ws.Range("J120").Value = 23249.75
$ws.Range("K120").Value = 35625
$ws.Range("L120").Value = 69749.25
$ws.Range("M120").Value = -30787
$ws.Range("N120").Value = -79425.25

$ws.Range("H135").Value = 1391.7222
$ws.Range("I135").Value = 1235.0834
$ws.Range("K135").Value = 11115.7506
$ws.Range("M135").Value = -8580.750599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5250
$ws.Range("J3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("N3").Value = -7232

$ws.Range("H74").Value = 50001
$ws.Range("J74").Value = 50001
$ws.Range("L74").Value = 50001
$ws.Range("N74").Value = -51873

$ws.Range("H77").Value = 50001
$ws.Range("J77").Value = 50001
$ws.Range("L77").Value = 150003
$ws.Range("N77").Value = -159363

$ws.Range("H97").Value = 1669.1666
$ws.Range("I97").Value = 1400.9
$ws.Range("J97").Value = 3010.5
$ws.Range("K97").Value = 1400.9
$ws.Range("L97").Value = 3010.5
$ws.Range("M97").Value = -904.9000000000001
$ws.Range("N97").Value = -4002.5

$ws.Range("H134").Value = 49399.8
$ws.Range("J134").Value = 49399.8
$ws.Range("L134").Value = 148199.4
$ws.Range("N134").Value = -153269.4

$ws.Range("H136").Value = 62662.5
$ws.Range("J136").Value = 62662.5
$ws.Range("L136").Value = 187987.5
$ws.Range("N136").Value = -193087.5

$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4263.857
$ws.Range("I7").Value = 3807.4167
$ws.Range("K7").Value = 3807.4167
$ws.Range("M7").Value = -3695.4167

$ws.Range("H16").Value = 1956
$ws.Range("J16").Value = 1800
$ws.Range("L16").Value = 1800
$ws.Range("N16").Value = -2140

$ws.Range("H61").Value = 4213.636
$ws.Range("I61").Value = 3837.25
$ws.Range("K61").Value = 3837.25
$ws.Range("M61").Value = -3635.25

$ws.Range("H93").Value = 19966
$ws.Range("I93").Value = 19965
$ws.Range("K93").Value = 19965
$ws.Range("M93").Value = -18717

$ws.Range("H113").Value = 4213.636
$ws.Range("I113").Value = 3837.25
$ws.Range("K113").Value = 3837.25
$ws.Range("M113").Value = -1667.25

$ws.Range("H126").Value = 4263.857
$ws.Range("I126").Value = 3807.4167
$ws.Range("K126").Value = 11422.2501
$ws.Range("M126").Value = -8952.250100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 77500
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 77500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 77500
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -77970

$ws.Range("H35").Value = 77500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 77500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 77500
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -78080

$ws.Range("H46").Value = 66700
$ws.Range("J46").Value = 66700
$ws.Range("L46").Value = 66700
$ws.Range("N46").Value = -67162

$ws.Range("H134").Value = 66700
$ws.Range("J134").Value = 66700
$ws.Range("L134").Value = 200100
$ws.Range("N134").Value = -205170
